# Refresh the crypto price/volume table with the latest scrape results.
# Column D sometimes holds numeric-looking text (e.g. "581.06", "1.00") that
# must stay literal text (matching the source inlineStr cells) instead of being
# parsed into a number, so those assignments use a leading apostrophe -
# exactly like typing `'581.06` into Excel - to force text entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '73.248.36'
$ws.Range('E2').Value = '  +2.04%  '

# Row 3
$ws.Range('D3').Value = '4.046.68'
$ws.Range('E3').Value = '  +1.22%  '

# Row 4
$ws.Range('E4').Value = '  -0.08%  '

# Row 5
$ws.Range('D5').Value = '''581.06'
$ws.Range('E5').Value = '  +9.96%  '

# Row 6
$ws.Range('D6').Value = '''151.56'
$ws.Range('E6').Value = '  +0.73%  '

# Row 7
$ws.Range('D7').Value = '4.042.20'
$ws.Range('E7').Value = '  +1.31%  '

# Row 8
$ws.Range('E8').Value = '  -0.01%  '

# Row 9
$ws.Range('D9').Value = '''0.999'
$ws.Range('E9').Value = '  -0.05%  '

# Row 10
$ws.Range('D10').Value = '''0.757'
$ws.Range('E10').Value = '  +2.05%  '

# Row 11
$ws.Range('E11').Value = '  -0.25%  '

# Row 12
$ws.Range('D12').Value = '''53.40'
$ws.Range('E12').Value = '  +12.65%  '

# Row 13
$ws.Range('E13').Value = '  -0.73%  '

# Row 14
$ws.Range('D14').Value = '''11.07'
$ws.Range('E14').Value = '  +4.46%  '

# Row 15
$ws.Range('D15').Value = '4.693.47'
$ws.Range('E15').Value = '  +1.13%  '

# Row 16
$ws.Range('D16').Value = '4.053.35'
$ws.Range('E16').Value = '  +1.38%  '

# Row 17
$ws.Range('D17').Value = '''14.28'
$ws.Range('E17').Value = '  +2.31%  '

# Row 18
$ws.Range('E18').Value = '  +3.96%  '

# Row 19
$ws.Range('E19').Value = '  +1.12%  '

# Row 20
$ws.Range('D20').Value = '73.174.96'
$ws.Range('E20').Value = '  +2.11%  '

# Row 21
$ws.Range('E21').Value = '  -0.35%  '

# Row 22
$ws.Range('D22').Value = '''440.90'
$ws.Range('E22').Value = '  +3.50%  '

# Row 23
$ws.Range('D23').Value = '''4.61'
$ws.Range('E23').Value = '  +10.89%  '

# Row 24
$ws.Range('D24').Value = '''97.66'
$ws.Range('E24').Value = '  +0.46%  '

# Row 25
$ws.Range('D25').Value = '''3.53'
$ws.Range('E25').Value = '  +1.99%  '

# Row 26
$ws.Range('E26').Value = '  +1.82%  '

# Row 27
$ws.Range('D27').Value = '''4.31'
$ws.Range('E27').Value = '  +20.36%  '

# Row 28
$ws.Range('D28').Value = '''11.54'
$ws.Range('E28').Value = '  +3.49%  '

# Row 29
$ws.Range('E29').Value = '  +2.74%  '

# Row 30
$ws.Range('D30').Value = '''5.97'
$ws.Range('E30').Value = '  +2.24%  '

# Row 31
$ws.Range('D31').Value = '''36.96'
$ws.Range('E31').Value = '  +1.28%  '

# Row 32
$ws.Range('D32').Value = '''7.97'
$ws.Range('E32').Value = '  +14.49%  '

# Row 33
$ws.Range('E33').Value = '  +4.33%  '

# Row 34
$ws.Range('D34').Value = '''13.65'
$ws.Range('E34').Value = '  +2.64%  '

# Row 35
$ws.Range('D35').Value = '''691.56'
$ws.Range('E35').Value = '  +2.40%  '

# Row 36
$ws.Range('D36').Value = '''48.51'
$ws.Range('E36').Value = '  +10.90%  '

# Row 37
$ws.Range('D37').Value = '''67.41'
$ws.Range('E37').Value = '  +3.28%  '

# Row 38
$ws.Range('E38').Value = '  +3.37%  '

# Row 39
$ws.Range('D39').Value = '0.0₃0886'
$ws.Range('E39').Value = '  +7.64%  '

# Row 40
$ws.Range('E40').Value = '  -1.30%  '

# Row 41
$ws.Range('D41').Value = '''11.27'
$ws.Range('E41').Value = '  +17.91%  '

# Row 42
$ws.Range('D42').Value = '''3.37'
$ws.Range('E42').Value = '  -0.77%  '

# Row 43
$ws.Range('D43').Value = '''1.00'
$ws.Range('E43').Value = '  +0.06%  '

# Row 44
$ws.Range('D44').Value = '''3.34'
$ws.Range('E44').Value = '  +6.27%  '

# Row 45
$ws.Range('D45').Value = '''0.0495'
$ws.Range('E45').Value = '  +2.35%  '

# Row 46
$ws.Range('E46').Value = '  +0.20%  '

# Row 47
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').Value = '''0.151'
$ws.Range('E47').Value = '  +1.00%  '

# Row 48
$ws.Range('B48').Value = 'Fetch.AI'
$ws.Range('C48').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D48').Value = '''2.76'
$ws.Range('E48').Value = '  +5.36%  '

# Row 49
$ws.Range('D49').Value = '''3.38'
$ws.Range('E49').Value = '  -0.67%  '

# Row 50
$ws.Range('D50').Value = '''3.50'
$ws.Range('E50').Value = '  +6.89%  '

# Row 51
$ws.Range('E51').Value = '  +2.97%  '
